# Update "Estado de Cuenta" workbook: keep only the first worker record and
# update the summary counters, then move the footer (signature) rows up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the worker detail rows that are no longer part of this statement
# (rows 17-23), leaving only the single remaining worker in row 16.
$ws.Range("B17:J23").EntireRow.Delete()

# Update "Valor Mora" total (now equal to the single remaining worker's value).
$ws.Range("E11").Value = 22533

# Update worker / period counters.
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1
